$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old row 6 values
$ws.Range("D2").Value = 44685
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 1500
$ws.Range("M2").Value = 1750
$ws.Range("P2").Value = 583

# Row 3 <- old row 4 values
$ws.Range("D3").Value = 44827
$ws.Range("J3").Value = 1200
$ws.Range("K3").Value = 2000
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = 2250
$ws.Range("P3").Value = 750

# Row 4 <- old row 2 values
$ws.Range("D4").Value = 44883
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 1800
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = 1900
$ws.Range("P4").Value = 633

# Row 6 <- old row 3 values
$ws.Range("D6").Value = 44848
$ws.Range("J6").Value = 1000
